$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header change
$ws.Range("G1").Value = "Journal"

# Row 2
$ws.Range("B2").Value = "Technical note: exploiting problem definition study for cyber security simulations"
$ws.Range("C2").Value = "Yilmaz Cankaya"
$ws.Range("D2").Value = "N/A"
$ws.Range("E2").Value = "10.1177/1548512915604585"
$ws.Range("F2").Value = "Restricted"
$ws.Range("G2").Value = "The Journal of Defense Modeling and Simulation: Applications, Methodology, Technology"

# Row 3
$ws.Range("B3").Value = "Moving beyond the sanctuary paradigm: Canada must face up to the reality of a contested and dangerous space environment"
$ws.Range("C3").Value = "Patrick Perron"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "10.1177/00207020231178394"
$ws.Range("F3").Value = "Restricted"
$ws.Range("G3").Value = "International Journal: Canada’s Journal of Global Policy Analysis"

# Row 4
$ws.Range("B4").Value = "Robust tracking strategy for nonlinear connected vehicle cyber-physical systems"
$ws.Range("C4").Value = "Yushi Yang, Meng Li, Yong Chen"
$ws.Range("D4").Value = "N/A"
$ws.Range("E4").Value = "10.1177/01423312231196642"
$ws.Range("F4").Value = "Restricted"
$ws.Range("G4").Value = "Transactions of the Institute of Measurement and Control"

# Row 5
$ws.Range("B5").Value = "Using network digital twins to improve cyber resilience of missions"
$ws.Range("C5").Value = "Rajive Bagrodia"
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "10.1177/15485129221131226"
$ws.Range("F5").Value = "Restricted"
$ws.Range("G5").Value = "The Journal of Defense Modeling and Simulation: Applications, Methodology, Technology"

# Row 6
$ws.Range("B6").Value = "Attrition rates and maneuver in agent-based simulation models"
$ws.Range("C6").Value = "David Ormrod, Benjamin Turnbull"
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "10.1177/1548512917692693"
$ws.Range("F6").Value = "Restricted"
$ws.Range("G6").Value = "The Journal of Defense Modeling and Simulation: Applications, Methodology, Technology"

# Row 7
$ws.Range("B7").Value = "A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"
$ws.Range("C7").Value = "Ramesh Sneka Nandhini, Ramanathan Lakshmanan"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "10.3233/JIFS-235809"
$ws.Range("F7").Value = "Restricted"
$ws.Range("G7").Value = "Journal of Intelligent & Fuzzy Systems: Applications in Engineering and Technology"

# Row 8
$ws.Range("B8").Value = "Accountability and cyber conflict: examining institutional constraints on the use of cyber proxies"
$ws.Range("C8").Value = "William Akoto"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "10.1177/07388942211051264"
$ws.Range("F8").Value = "Restricted"
$ws.Range("G8").Value = "Conflict Management and Peace Science"

# Row 9
$ws.Range("B9").Value = "Wargaming the use of intermediate force capabilities in the gray zone"
$ws.Range("C9").Value = "Kyle D Christensen, Peter Dobias"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "10.1177/15485129211010227"
$ws.Range("F9").Value = "Restricted"
$ws.Range("G9").Value = "The Journal of Defense Modeling and Simulation: Applications, Methodology, Technology"

# Row 10
$ws.Range("B10").Value = "Cyber and contentious politics: Evidence from the US radical environmental movement"
$ws.Range("C10").Value = "Thomas Zeitzoff, Grace Gold"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "10.1177/00223433231221426"
$ws.Range("F10").Value = "Restricted"
$ws.Range("G10").Value = "Journal of Peace Research"

# Row 11
$ws.Range("B11").Value = "Internet of Things, cybersecurity and governing wicked problems: learning from climate change governance"
$ws.Range("C11").Value = "Madeline Carr, Feja Lesniewska"
$ws.Range("D11").Value = "N/A"
$ws.Range("E11").Value = "10.1177/0047117820948247"
$ws.Range("F11").Value = "Open Access"
$ws.Range("G11").Value = "International Relations"
